$wb = $excel.ActiveWorkbook

# --- Add the new "OrderPayload" sheet at the end of the tab strip ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$payloadSheet = $wb.Worksheets.Add($null, $lastSheet)
$payloadSheet.Name = "OrderPayload"

$payloadSheet.Range("C5").Value = "Order Payload"
$payloadSheet.Range("D5").Value = "onSubmit"

$payloadSheet.Columns.Item(3).ColumnWidth = 16
$payloadSheet.Columns.Item(4).ColumnWidth = 13.666666666666666

$payloadSheet.Range("D6").Select() | Out-Null

# --- Extend the "Sheet3" navigation list with the new menu entries ---
$navSheet = $wb.Worksheets.Item("Sheet3")

$navSheet.Range("C13").Value = "Management"
$navSheet.Range("C14").Value = "Dashboard"
$navSheet.Range("C15").Value = "Documents"
$navSheet.Range("C16").Value = "Products"
$navSheet.Range("C17").Value = "Stock"
$navSheet.Range("C18").Value = "Reposrting"
$navSheet.Range("C19").Value = "Customers & suppliers"
$navSheet.Range("C20").Value = "Promotions & actions"
$navSheet.Range("C21").Value = "Users & security"
$navSheet.Range("C22").Value = "Cash Registers"

$navSheet.Range("C14").Select() | Out-Null
$navSheet.Activate() | Out-Null
